$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.753.18'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.438.97'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.600'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('D9').Value = '3.435.30'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.131'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.408'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('D13').Value = '4.031.53'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.29%  '
$ws.Range('D16').Value = '65.787.15'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000170'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '3.436.14'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.73'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '366.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.531'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000121'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.177'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.56%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '23.59'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.95%  '
$ws.Range('E35').Value = '  -2.41%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.878'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '28.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.62%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.59'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = '2.757.22'
$ws.Range('E42').Value = '  +3.05%  '
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0677'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0288'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '325.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.71%  '
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.80%  '
